$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceOne = 1

# 1. Update activation date
$d.Content.Find.Execute(
    "Ativação: 01/01/2022", $true, $false, $false, $false, $false, $true,
    $wdFindContinue, $false, "Ativação: 01/01/2024", $wdReplaceOne)

# 2. Add two docentes around the existing one (Carlos Angelo Nunes before,
#    Luiz Tadeu Fernandes Eleno after), each former line ending in a manual
#    line break except the last added line.
$br = [char]11
$r = $d.Content
$r.Find.Execute(
    "5009972 - Gilberto Carvalho Coelho", $true, $false, $false, $false,
    $false, $true, $wdFindContinue, $false, "", 0)
$before = $r.Duplicate
$before.Collapse(1)
$before.InsertBefore("3577649 - Carlos Angelo Nunes" + $br)

$r2 = $d.Content
$r2.Find.Execute(
    "5009972 - Gilberto Carvalho Coelho", $true, $false, $false, $false,
    $false, $true, $wdFindContinue, $false, "", 0)
$after = $r2.Duplicate
$after.Collapse(0)
$after.InsertAfter($br + "1176388 - Luiz Tadeu Fernandes Eleno")

# 3. Programa resumido: add spaces and a new item F
$d.Content.Find.Execute(
    "A. Introdução; teoria básica de equilíbrio de fases;B. Sistemas unários;C. Sistemas binários;D. Sistemas ternários;E. Cálculo termodinâmico de diagramas de fases.",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "A. Introdução; teoria básica de equilíbrio de fases;B. Sistemas unários; C. Sistemas binários; D. Sistemas ternários; E. Cálculo termodinâmico de diagramas de fases; F. Trabalho Prático.",
    $wdReplaceOne)

# 4a. Programa: remove the space after "invariantes;"
$d.Content.Find.Execute(
    "invariantes; 3. Sistemas binários", $true, $false, $false, $false,
    $false, $true, $wdFindContinue, $false,
    "invariantes;3. Sistemas binários", $wdReplaceOne)

# 4b. Programa: append "Trabalho prático." at the end
$d.Content.Find.Execute(
    "12. Cálculo termodinâmico de diagramas de fases.", $true, $false,
    $false, $false, $false, $true, $wdFindContinue, $false,
    "12. Cálculo termodinâmico de diagramas de fases; Trabalho prático.",
    $wdReplaceOne)

# 5. Método text
$d.Content.Find.Execute(
    "O curso será ministrado na forma de aulas expositivas e aulas práticas em laboratório envolvendo preparação de amostras e caracterização microestrutural. Os resultados das aulas práticas serão apresentados oralmente e sujeitos a avaliação (T).",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "O curso será ministrado na forma de aulas expositivas e aulas práticas em laboratório envolvendo preparação de amostras e caracterização microestrutural. Os resultados das aulas práticas serão apresentados oralmente e por escrito. Questionários e listas de exercícios serão elaborados para serem respondidos individualmente ou em grupo. Avaliações escritas serão realizadas para resolução individual.",
    $wdReplaceOne)

# 6. Critério text
$d.Content.Find.Execute(
    "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF) juntamente com a avaliação do trabalho prático (T). O critério para a nota final é:NF=((P1*0,8)+(T*0,2)+P2*1)/2",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "As avaliações individuais, a participação nas resoluções dos exercícios e repostas aos questionários assim como a condução do trabalho prático e a apresentação dos resultados nas formas oral e escrita serão agrupadas em duas notas (N1 e N2) que comporão a nota final (NF). O critério para cálculo da nota final é: NF = (N1+ N2)/2Serão aprovados os alunos com NF ≥ 5,0Serão reprovados os alunos com NF < 3,0",
    $wdReplaceOne)

# 7. Norma de recuperação text
$d.Content.Find.Execute(
    "Para os alunos que obtiverem 3,0≤NF<5,0, será aplicada uma avaliação de recuperação (R) que levará ao cálculo da média final (MF) com o seguinte critério:MF=(NF+R)/2",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Será aplicada recuperação para os alunos que obtiverem NF entre 3,0 e 4,9. A nota pós recuperação será calculada pela média aritmética com a nota final NF.",
    $wdReplaceOne)
